$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-11-02 Sunday"; new="2025-11-03 Monday"},
    @{old="99÷6=16, 3"; new="63÷8=7, 7"},
    @{old="12÷2=6, 0"; new="28÷4=7, 0"},
    @{old="16÷6=2, 4"; new="57÷6=9, 3"},
    @{old="99÷4=24, 3"; new="95÷5=19, 0"},
    @{old="38÷2=19, 0"; new="36÷7=5, 1"},
    @{old="72÷4=18, 0"; new="66÷2=33, 0"},
    @{old="10÷7=1, 3"; new="76÷5=15, 1"},
    @{old="51÷3=17, 0"; new="93÷5=18, 3"},
    @{old="84÷5=16, 4"; new="84÷8=10, 4"},
    @{old="56÷8=7, 0"; new="38÷8=4, 6"},
    @{old="17÷9=1, 8"; new="34÷5=6, 4"},
    @{old="87÷8=10, 7"; new="41÷9=4, 5"},
    @{old="51÷8=6, 3"; new="76÷4=19, 0"},
    @{old="25÷4=6, 1"; new="72÷7=10, 2"},
    @{old="87÷5=17, 2"; new="50÷9=5, 5"},
    @{old="54÷2=27, 0"; new="60÷3=20, 0"},
    @{old="23÷8=2, 7"; new="15÷5=3, 0"},
    @{old="33÷8=4, 1"; new="43÷5=8, 3"},
    @{old="46÷7=6, 4"; new="30÷4=7, 2"},
    @{old="44÷4=11, 0"; new="19÷2=9, 1"},
    @{old="29÷6=4, 5"; new="12÷6=2, 0"},
    @{old="83÷5=16, 3"; new="33÷5=6, 3"},
    @{old="81÷3=27, 0"; new="18÷6=3, 0"},
    @{old="57÷7=8, 1"; new="70÷8=8, 6"},
    @{old="30÷7=4, 2"; new="26÷9=2, 8"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
